# "teste" -> "Teste teste", laid out as three separate runs:
#   "T" | "este" | " teste"
#
# Strategy: perform the plain text edit first (capitalize the first
# letter, append " teste"), which Word naturally keeps as one run since
# the formatting never changes. Then force run boundaries exactly where
# the target runs must split by dropping a (temporary) bookmark at each
# boundary and immediately deleting it -- the bookmark's insertion point
# splits the underlying run, and the split survives the bookmark's
# removal.

$d = $word.ActiveDocument

# 1) Capitalize the leading "t" -> "T".
$first = $d.Range(0, 1)
$first.Text = "T"

# 2) Append " teste" right after the original word (now "Teste").
$tail = $d.Range(5, 5)
$tail.InsertBefore(" teste")

# 3) Split the single run into "T" | "este" | " teste" by planting and
#    removing bookmarks at the two boundaries (offsets 1 and 5).
$d.Bookmarks.Add("_split1", $d.Range(1, 1))
$d.Bookmarks.Add("_split2", $d.Range(5, 5))
$d.Bookmarks("_split1").Delete()
$d.Bookmarks("_split2").Delete()
